# Update summary counts on the "Inscricoes" worksheet to reflect newly
# registered/processed applications (Commit via gitrun.py em 2024-09-21 12:00:45)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Processo "Iluminacao Publica" / Nivel Superior line): Inscritos 16 -> 17
$ws.Range("E4").Value = 17

# Row 12: Inscritos 23 -> 24, Pagos 7 -> 8, Inscrições homologadas 7 -> 8
$ws.Range("E12").Value = 24
$ws.Range("F12").Value = 8
$ws.Range("H12").Value = 8

# Row 16: Inscritos 294 -> 296
$ws.Range("E16").Value = 296
